# Auto-generated: refresh Market Board derived profit figures across all class sheets
# (scheduled runner sync of currentAveragePrice / Leve price / profit columns)
$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Cells.Item(76, 8).Value = 2996.6667  # H76: was 3165.5293
$ws_ALC.Cells.Item(76, 9).Value = 2995.5  # I76: was 3003
$ws_ALC.Cells.Item(76, 10).Value = 3003.2  # J76: was 3187.2
$ws_ALC.Cells.Item(76, 11).Value = 2995.5  # K76: was 3003
$ws_ALC.Cells.Item(76, 12).Value = 3003.2  # L76: was 3187.2
$ws_ALC.Cells.Item(76, 13).Value = -2680.5  # M76: was -2688
$ws_ALC.Cells.Item(76, 14).Value = -3633.2  # N76: was -3817.2
$ws_ALC.Cells.Item(79, 8).Value = 2996.6667  # H79: was 3165.5293
$ws_ALC.Cells.Item(79, 9).Value = 2995.5  # I79: was 3003
$ws_ALC.Cells.Item(79, 10).Value = 3003.2  # J79: was 3187.2
$ws_ALC.Cells.Item(79, 11).Value = 2995.5  # K79: was 3003
$ws_ALC.Cells.Item(79, 12).Value = 3003.2  # L79: was 3187.2
$ws_ALC.Cells.Item(79, 13).Value = -1903.5  # M79: was -1911
$ws_ALC.Cells.Item(79, 14).Value = -5187.2  # N79: was -5371.2
$ws_ALC.Cells.Item(137, 8).Value = 2082954.2  # H137: was 1976223.5
$ws_ALC.Cells.Item(137, 9).Value = 5495406  # I137: was 5495493
$ws_ALC.Cells.Item(137, 10).Value = 5809.609  # J137: was 5432.56
$ws_ALC.Cells.Item(137, 11).Value = 16486218  # K137: was 16486479
$ws_ALC.Cells.Item(137, 12).Value = 17428.827  # L137: was 16297.68
$ws_ALC.Cells.Item(137, 13).Value = -16483668  # M137: was -16483929
$ws_ALC.Cells.Item(137, 14).Value = -22528.827  # N137: was -21397.68

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Cells.Item(63, 8).Value = 3168.238  # H63: was 3034.7083
$ws_ARM.Cells.Item(63, 9).Value = 2257.4614  # I63: was 2280.5386
$ws_ARM.Cells.Item(63, 10).Value = 4648.25  # J63: was 3926
$ws_ARM.Cells.Item(63, 11).Value = 2257.4614  # K63: was 2280.5386
$ws_ARM.Cells.Item(63, 12).Value = 4648.25  # L63: was 3926
$ws_ARM.Cells.Item(63, 13).Value = -1571.4614  # M63: was -1594.5386
$ws_ARM.Cells.Item(63, 14).Value = -6020.25  # N63: was -5298
$ws_ARM.Cells.Item(66, 8).Value = 3168.238  # H66: was 3034.7083
$ws_ARM.Cells.Item(66, 9).Value = 2257.4614  # I66: was 2280.5386
$ws_ARM.Cells.Item(66, 10).Value = 4648.25  # J66: was 3926
$ws_ARM.Cells.Item(66, 11).Value = 11287.307  # K66: was 11402.693
$ws_ARM.Cells.Item(66, 12).Value = 23241.25  # L66: was 19630
$ws_ARM.Cells.Item(66, 13).Value = -7855.307000000001  # M66: was -7970.692999999999
$ws_ARM.Cells.Item(66, 14).Value = -30105.25  # N66: was -26494
$ws_ARM.Cells.Item(97, 8).Value = 2655.5  # H97: was 1042.5217
$ws_ARM.Cells.Item(97, 9).Value = 2300  # I97: was 953.0454999999999
$ws_ARM.Cells.Item(97, 10).Value = 3011  # J97: was 3011
$ws_ARM.Cells.Item(97, 11).Value = 2300  # K97: was 953.0454999999999
$ws_ARM.Cells.Item(97, 12).Value = 3011  # L97: was 3011
$ws_ARM.Cells.Item(97, 13).Value = -1804  # M97: was -457.0454999999999
$ws_ARM.Cells.Item(97, 14).Value = -4003  # N97: was -4003
$ws_ARM.Cells.Item(132, 8).Value = 17243668  # H132: was 16668987
$ws_ARM.Cells.Item(132, 9).Value = 26317142  # I132: was 23810934
$ws_ARM.Cells.Item(132, 10).Value = 4069.3  # J132: was 4443.6665
$ws_ARM.Cells.Item(132, 11).Value = 78951426  # K132: was 71432802
$ws_ARM.Cells.Item(132, 12).Value = 12207.9  # L132: was 13330.9995
$ws_ARM.Cells.Item(132, 13).Value = -78948896  # M132: was -71430272
$ws_ARM.Cells.Item(132, 14).Value = -17267.9  # N132: was -18390.9995

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Cells.Item(20, 8).Value = 2238.7222  # H20: was 2282.1765
$ws_BSM.Cells.Item(20, 9).Value = 1666.8572  # I20: was 1694.6666
$ws_BSM.Cells.Item(20, 10).Value = 2602.6365  # J20: was 2602.6365
$ws_BSM.Cells.Item(20, 11).Value = 1666.8572  # K20: was 1694.6666
$ws_BSM.Cells.Item(20, 12).Value = 2602.6365  # L20: was 2602.6365
$ws_BSM.Cells.Item(20, 13).Value = -1419.8572  # M20: was -1447.6666
$ws_BSM.Cells.Item(20, 14).Value = -3096.6365  # N20: was -3096.6365
$ws_BSM.Cells.Item(86, 8).Value = 2255.2856  # H86: was 2214.2
$ws_BSM.Cells.Item(86, 9).Value = 2216  # I86: was 2119.2856
$ws_BSM.Cells.Item(86, 10).Value = 2353.5  # J86: was 2435.6667
$ws_BSM.Cells.Item(86, 11).Value = 2216  # K86: was 2119.2856
$ws_BSM.Cells.Item(86, 12).Value = 2353.5  # L86: was 2435.6667
$ws_BSM.Cells.Item(86, 13).Value = -1093  # M86: was -996.2856000000002
$ws_BSM.Cells.Item(86, 14).Value = -4599.5  # N86: was -4681.6667
$ws_BSM.Cells.Item(89, 8).Value = 2255.2856  # H89: was 2214.2
$ws_BSM.Cells.Item(89, 9).Value = 2216  # I89: was 2119.2856
$ws_BSM.Cells.Item(89, 10).Value = 2353.5  # J89: was 2435.6667
$ws_BSM.Cells.Item(89, 11).Value = 11080  # K89: was 10596.428
$ws_BSM.Cells.Item(89, 12).Value = 11767.5  # L89: was 12178.3335
$ws_BSM.Cells.Item(89, 13).Value = -5464  # M89: was -4980.428
$ws_BSM.Cells.Item(89, 14).Value = -22999.5  # N89: was -23410.3335
$ws_BSM.Cells.Item(134, 8).Value = 2865.7546  # H134: was 2576.1428
$ws_BSM.Cells.Item(134, 9).Value = 1806.381  # I134: was 1568.4849
$ws_BSM.Cells.Item(134, 10).Value = 3560.9688  # J134: was 3684.5667
$ws_BSM.Cells.Item(134, 11).Value = 5419.143  # K134: was 4705.4547
$ws_BSM.Cells.Item(134, 12).Value = 10682.9064  # L134: was 11053.7001
$ws_BSM.Cells.Item(134, 13).Value = -2884.143  # M134: was -2170.4547
$ws_BSM.Cells.Item(134, 14).Value = -15752.9064  # N134: was -16123.7001

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Cells.Item(31, 8).Value = 5160.27  # H31: was 5306.3276
$ws_CRP.Cells.Item(31, 9).Value = 2214.2307  # I31: was 2627.4
$ws_CRP.Cells.Item(31, 10).Value = 5926.24  # J31: was 5831.608
$ws_CRP.Cells.Item(31, 11).Value = 2214.2307  # K31: was 2627.4
$ws_CRP.Cells.Item(31, 12).Value = 5926.24  # L31: was 5831.608
$ws_CRP.Cells.Item(31, 13).Value = -1919.2307  # M31: was -2332.4
$ws_CRP.Cells.Item(31, 14).Value = -6516.24  # N31: was -6421.608
$ws_CRP.Cells.Item(34, 8).Value = 5160.27  # H34: was 5306.3276
$ws_CRP.Cells.Item(34, 9).Value = 2214.2307  # I34: was 2627.4
$ws_CRP.Cells.Item(34, 10).Value = 5926.24  # J34: was 5831.608
$ws_CRP.Cells.Item(34, 11).Value = 2214.2307  # K34: was 2627.4
$ws_CRP.Cells.Item(34, 12).Value = 5926.24  # L34: was 5831.608
$ws_CRP.Cells.Item(34, 13).Value = -2012.2307  # M34: was -2425.4
$ws_CRP.Cells.Item(34, 14).Value = -6330.24  # N34: was -6235.608
$ws_CRP.Cells.Item(134, 8).Value = 420176.38  # H134: was 555203.3
$ws_CRP.Cells.Item(134, 9).Value = 441638.44  # I134: was 614378.0600000001
$ws_CRP.Cells.Item(134, 10).Value = 282819.2  # J134: was 282999.2
$ws_CRP.Cells.Item(134, 11).Value = 1324915.32  # K134: was 1843134.18
$ws_CRP.Cells.Item(134, 12).Value = 848457.6000000001  # L134: was 848997.6000000001
$ws_CRP.Cells.Item(134, 13).Value = -1322380.32  # M134: was -1840599.18
$ws_CRP.Cells.Item(134, 14).Value = -853527.6000000001  # N134: was -854067.6000000001

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Cells.Item(69, 8).Value = 68629790  # H69: was 73531656
$ws_CUL.Cells.Item(69, 9).Value = 0  # I69: was 0
$ws_CUL.Cells.Item(69, 10).Value = 68629790  # J69: was 73531656
$ws_CUL.Cells.Item(69, 11).Value = 0  # K69: was 0
$ws_CUL.Cells.Item(69, 12).Value = 205889370  # L69: was 220594968
$ws_CUL.Cells.Item(69, 14).Value = -205890992  # N69: was -220596590
$ws_CUL.Cells.Item(72, 8).Value = 68629790  # H72: was 73531656
$ws_CUL.Cells.Item(72, 9).Value = 0  # I72: was 0
$ws_CUL.Cells.Item(72, 10).Value = 68629790  # J72: was 73531656
$ws_CUL.Cells.Item(72, 11).Value = 0  # K72: was 0
$ws_CUL.Cells.Item(72, 12).Value = 617668110  # L72: was 661784904
$ws_CUL.Cells.Item(72, 14).Value = -617676222  # N72: was -661793016
$ws_CUL.Cells.Item(82, 8).Value = 2000  # H82: was 3000
$ws_CUL.Cells.Item(82, 9).Value = 1000  # I82: was 0
$ws_CUL.Cells.Item(82, 10).Value = 3000  # J82: was 3000
$ws_CUL.Cells.Item(82, 11).Value = 3000  # K82: was 0
$ws_CUL.Cells.Item(82, 12).Value = 9000  # L82: was 9000
$ws_CUL.Cells.Item(82, 13).Value = -2594  # M82: was None
$ws_CUL.Cells.Item(82, 14).Value = -9812  # N82: was -9812
$ws_CUL.Cells.Item(85, 8).Value = 2000  # H85: was 3000
$ws_CUL.Cells.Item(85, 9).Value = 1000  # I85: was 0
$ws_CUL.Cells.Item(85, 10).Value = 3000  # J85: was 3000
$ws_CUL.Cells.Item(85, 11).Value = 3000  # K85: was 0
$ws_CUL.Cells.Item(85, 12).Value = 9000  # L85: was 9000
$ws_CUL.Cells.Item(85, 13).Value = -1596  # M85: was None
$ws_CUL.Cells.Item(85, 14).Value = -11808  # N85: was -11808
$ws_CUL.Cells.Item(113, 8).Value = 6692.1177  # H113: was 5517.4287
$ws_CUL.Cells.Item(113, 9).Value = 50701  # I113: was 14857.429
$ws_CUL.Cells.Item(113, 10).Value = 824.26666  # J113: was 847.4286
$ws_CUL.Cells.Item(113, 11).Value = 152103  # K113: was 44572.287
$ws_CUL.Cells.Item(113, 12).Value = 2472.79998  # L113: was 2542.2858
$ws_CUL.Cells.Item(113, 13).Value = -149933  # M113: was -42402.287
$ws_CUL.Cells.Item(113, 14).Value = -6812.79998  # N113: was -6882.2858
$ws_CUL.Cells.Item(115, 8).Value = 4077.9614  # H115: was 4061.08
$ws_CUL.Cells.Item(115, 9).Value = 0  # I115: was 1400
$ws_CUL.Cells.Item(115, 10).Value = 4077.9614  # J115: was 4171.9585
$ws_CUL.Cells.Item(115, 11).Value = 0  # K115: was 4200
$ws_CUL.Cells.Item(115, 12).Value = 12233.8842  # L115: was 12515.8755
$ws_CUL.Cells.Item(115, 13).ClearContents()  # M115: was -3025
$ws_CUL.Cells.Item(115, 14).Value = -14583.8842  # N115: was -14865.8755
$ws_CUL.Cells.Item(121, 8).Value = 97741.81  # H121: was 92116.32000000001
$ws_CUL.Cells.Item(121, 9).Value = 333.33334  # I121: was 610
$ws_CUL.Cells.Item(121, 10).Value = 107818.555  # J121: was 100971.77
$ws_CUL.Cells.Item(121, 11).Value = 1000.00002  # K121: was 1830
$ws_CUL.Cells.Item(121, 12).Value = 323455.665  # L121: was 302915.31
$ws_CUL.Cells.Item(121, 13).Value = 309.9999799999999  # M121: was -520
$ws_CUL.Cells.Item(121, 14).Value = -326075.665  # N121: was -305535.31

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Cells.Item(70, 8).Value = 5145.1816  # H70: was 5172.591
$ws_GSM.Cells.Item(70, 9).Value = 5171.143  # I70: was 5199.857
$ws_GSM.Cells.Item(70, 10).Value = 4600  # J70: was 4600
$ws_GSM.Cells.Item(70, 11).Value = 5171.143  # K70: was 5199.857
$ws_GSM.Cells.Item(70, 12).Value = 4600  # L70: was 4600
$ws_GSM.Cells.Item(70, 13).Value = -4901.143  # M70: was -4929.857
$ws_GSM.Cells.Item(70, 14).Value = -5140  # N70: was -5140
$ws_GSM.Cells.Item(73, 8).Value = 5145.1816  # H73: was 5172.591
$ws_GSM.Cells.Item(73, 9).Value = 5171.143  # I73: was 5199.857
$ws_GSM.Cells.Item(73, 10).Value = 4600  # J73: was 4600
$ws_GSM.Cells.Item(73, 11).Value = 5171.143  # K73: was 5199.857
$ws_GSM.Cells.Item(73, 12).Value = 4600  # L73: was 4600
$ws_GSM.Cells.Item(73, 13).Value = -4235.143  # M73: was -4263.857
$ws_GSM.Cells.Item(73, 14).Value = -6472  # N73: was -6472
$ws_GSM.Cells.Item(102, 8).Value = 823.75  # H102: was 800.6667
$ws_GSM.Cells.Item(102, 9).Value = 832.4  # I102: was 760.8
$ws_GSM.Cells.Item(102, 10).Value = 809.3333  # J102: was 1000
$ws_GSM.Cells.Item(102, 11).Value = 832.4  # K102: was 760.8
$ws_GSM.Cells.Item(102, 12).Value = 809.3333  # L102: was 1000
$ws_GSM.Cells.Item(102, 13).Value = 789.6  # M102: was 861.2
$ws_GSM.Cells.Item(102, 14).Value = -4053.3333  # N102: was -4244
$ws_GSM.Cells.Item(132, 8).Value = 23259546  # H132: was 22225818
$ws_GSM.Cells.Item(132, 9).Value = 37040516  # I132: was 35717650
$ws_GSM.Cells.Item(132, 10).Value = 4156.3125  # J132: was 3978.5293
$ws_GSM.Cells.Item(132, 11).Value = 111121548  # K132: was 107152950
$ws_GSM.Cells.Item(132, 12).Value = 12468.9375  # L132: was 11935.5879
$ws_GSM.Cells.Item(132, 13).Value = -111119018  # M132: was -107150420
$ws_GSM.Cells.Item(132, 14).Value = -17528.9375  # N132: was -16995.5879

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Cells.Item(24, 8).Value = 9906  # H24: was 9900
$ws_LTW.Cells.Item(24, 9).Value = 9906  # I24: was 9900
$ws_LTW.Cells.Item(24, 10).Value = 0  # J24: was 0
$ws_LTW.Cells.Item(24, 11).Value = 9906  # K24: was 9900
$ws_LTW.Cells.Item(24, 12).Value = 0  # L24: was 0
$ws_LTW.Cells.Item(24, 13).Value = -9563  # M24: was -9557
$ws_LTW.Cells.Item(82, 8).Value = 4630689.5  # H82: was 5556755.5
$ws_LTW.Cells.Item(82, 9).Value = 1007  # I82: was 1124.6666
$ws_LTW.Cells.Item(82, 10).Value = 20834580  # J82: was 27779278
$ws_LTW.Cells.Item(82, 11).Value = 1007  # K82: was 1124.6666
$ws_LTW.Cells.Item(82, 12).Value = 20834580  # L82: was 27779278
$ws_LTW.Cells.Item(82, 13).Value = -646  # M82: was -763.6666
$ws_LTW.Cells.Item(82, 14).Value = -20835302  # N82: was -27780000
$ws_LTW.Cells.Item(85, 8).Value = 4630689.5  # H85: was 5556755.5
$ws_LTW.Cells.Item(85, 9).Value = 1007  # I85: was 1124.6666
$ws_LTW.Cells.Item(85, 10).Value = 20834580  # J85: was 27779278
$ws_LTW.Cells.Item(85, 11).Value = 1007  # K85: was 1124.6666
$ws_LTW.Cells.Item(85, 12).Value = 20834580  # L85: was 27779278
$ws_LTW.Cells.Item(85, 13).Value = 241  # M85: was 123.3334
$ws_LTW.Cells.Item(85, 14).Value = -20837076  # N85: was -27781774
$ws_LTW.Cells.Item(93, 8).Value = 1782.091  # H93: was 2063.5715
$ws_LTW.Cells.Item(93, 9).Value = 999  # I93: was 2330
$ws_LTW.Cells.Item(93, 10).Value = 1860.4  # J93: was 1990.909
$ws_LTW.Cells.Item(93, 11).Value = 999  # K93: was 2330
$ws_LTW.Cells.Item(93, 12).Value = 1860.4  # L93: was 1990.909
$ws_LTW.Cells.Item(93, 13).Value = 249  # M93: was -1082
$ws_LTW.Cells.Item(93, 14).Value = -4356.4  # N93: was -4486.909
$ws_LTW.Cells.Item(136, 8).Value = 1465.4255  # H136: was 1694.4054
$ws_LTW.Cells.Item(136, 9).Value = 1088.9474  # I136: was 1248.3793
$ws_LTW.Cells.Item(136, 10).Value = 3055  # J136: was 3311.25
$ws_LTW.Cells.Item(136, 11).Value = 3266.8422  # K136: was 3745.1379
$ws_LTW.Cells.Item(136, 12).Value = 9165  # L136: was 9933.75
$ws_LTW.Cells.Item(136, 13).Value = -716.8422  # M136: was -1195.1379
$ws_LTW.Cells.Item(136, 14).Value = -14265  # N136: was -15033.75

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Cells.Item(5, 8).Value = 4800000  # H5: was 4600000
$ws_WVR.Cells.Item(5, 9).Value = 6000000  # I5: was 5500000
$ws_WVR.Cells.Item(5, 10).Value = 4000000  # J5: was 4000000
$ws_WVR.Cells.Item(5, 11).Value = 6000000  # K5: was 5500000
$ws_WVR.Cells.Item(5, 12).Value = 4000000  # L5: was 4000000
$ws_WVR.Cells.Item(5, 13).Value = -5999888  # M5: was -5499888
$ws_WVR.Cells.Item(5, 14).Value = -4000224  # N5: was -4000224
$ws_WVR.Cells.Item(96, 8).Value = 1200  # H96: was 0
$ws_WVR.Cells.Item(96, 9).Value = 0  # I96: was 0
$ws_WVR.Cells.Item(96, 10).Value = 1200  # J96: was 0
$ws_WVR.Cells.Item(96, 11).Value = 0  # K96: was 0
$ws_WVR.Cells.Item(96, 12).Value = 1200  # L96: was 0
$ws_WVR.Cells.Item(96, 14).Value = -3946  # N96: was None
$ws_WVR.Cells.Item(113, 8).Value = 1136  # H113: was 1181.3334
$ws_WVR.Cells.Item(113, 9).Value = 1222  # I113: was 1222
$ws_WVR.Cells.Item(113, 10).Value = 1050  # J113: was 1100
$ws_WVR.Cells.Item(113, 11).Value = 3666  # K113: was 3666
$ws_WVR.Cells.Item(113, 12).Value = 3150  # L113: was 3300
$ws_WVR.Cells.Item(113, 13).Value = -1496  # M113: was -1496
$ws_WVR.Cells.Item(113, 14).Value = -7490  # N113: was -7640
$ws_WVR.Cells.Item(122, 8).Value = 966.6667  # H122: was 1000
$ws_WVR.Cells.Item(122, 9).Value = 950  # I122: was 1000
$ws_WVR.Cells.Item(122, 10).Value = 1000  # J122: was 0
$ws_WVR.Cells.Item(122, 11).Value = 2850  # K122: was 3000
$ws_WVR.Cells.Item(122, 12).Value = 3000  # L122: was 0
$ws_WVR.Cells.Item(122, 13).Value = -400  # M122: was -550
$ws_WVR.Cells.Item(122, 14).Value = -7900  # N122: was None

